# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets.
# Row -> (old, new) values: 2:211->216, 4:12866->12892, 5:1332->1334, 6:190->194,
# 10:217->218, 11:465->466, 16:396->401, 17:5489->5501, 19:42->47, 21:29->30, 23:120->124

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 216
    4  = 12892
    5  = 1334
    6  = 194
    10 = 218
    11 = 466
    16 = 401
    17 = 5501
    19 = 47
    21 = 30
    23 = 124
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
